$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header row cells: "_old" suffix -> "_FV2404", "_new" suffix -> "_FV2410"
for ($col = 1; $col -le 10; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = ($cell.Value2 -replace '_old$', '_FV2404')
}
for ($col = 12; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = ($cell.Value2 -replace '_new$', '_FV2410')
}

# Turn the used range into an Excel Table (ListObject) so the header names
# become the table's column headers.
$tbl = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $ws.Range("A1:U66"),
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$tbl.Name = "Table1"

# Freeze the header row (row 1) so it stays visible while scrolling.
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
